# Wed, May 13, 2020  5:07:59 PM
#
# 1) Re-theme the deck: the slide master's theme ("Integral" / Red Violet
#    colour scheme) is replaced with the stock "Office Theme" colour
#    scheme (font scheme and format scheme are identical between the two
#    themes already, so only the 10 colours that actually differ need to
#    move - dk1/lt1 are shared by both palettes).
# 2) Three tables (slides 14, 15, 16) switch from the deck's custom
#    table style to the built-in table style
#    {D5D7D03D-AC91-4754-BCE0-CEDB40149CC1}.

$p = $ppt.ActivePresentation

# --- 1. Swap the active theme's colour scheme to "Office" -----------------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

$officeColors = @{
    1  = 0x000000   # dk1
    2  = 0xFFFFFF   # lt1
    3  = 0x6A5444   # dk2      (44546A, stored BGR)
    4  = 0xE6E6E7   # lt2      (E7E6E6)
    5  = 0xD59B5B   # accent1  (5B9BD5)
    6  = 0x317DED   # accent2  (ED7D31)
    7  = 0xA5A5A5   # accent3  (A5A5A5)
    8  = 0x00C0FF   # accent4  (FFC000)
    9  = 0xC47244   # accent5  (4472C4)
    10 = 0x47AD70   # accent6  (70AD47)
    11 = 0xC16305   # hlink    (0563C1)
    12 = 0x724F95   # folHlink (954F72)
}

foreach ($idx in 1..12) {
    $colorScheme.Item($idx).RGB = $officeColors[$idx]
}

# --- 2. Point the three tables at the built-in table style ----------------
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shapeIdx in 1..$slide.Shapes.Count) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{D5D7D03D-AC91-4754-BCE0-CEDB40149CC1}")
        }
    }
}
